$wb = $excel.ActiveWorkbook

# --- Create the new "R3_P7" sheet as a copy of "R3_P6" (placed after it) ---
$srcSheet = $wb.Worksheets.Item("R3_P6")
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "R3_P7"

# --- Update the header/info block on the new sheet ---
$newSheet.Range("B1").Value = "2018-07-01"
$newSheet.Range("B2").Value = "DSPR Cohorts"
$newSheet.Range("B4").Value = "R3_P7"
$newSheet.Range("B5").Value = "vials with DSPR lines per cohort"
$newSheet.Range("B6").Value = "SURF NB#001 pg 92"

# --- Replace the sample-grid contents with the new DSPR cohort labels ---
$newSheet.Range("C9").Value = "C1A7"
$newSheet.Range("D9").Value = "C1A2"
$newSheet.Range("E9").Value = "C1B3"
$newSheet.Range("F9").Value = "C1A5"
$newSheet.Range("G9").Value = "C1A6"
$newSheet.Range("H9").Value = "C1A4"
$newSheet.Range("I9").Value = "C1A1"
$newSheet.Range("J9").Value = "C1B6"
$newSheet.Range("K9").Value = "C1B4"

$newSheet.Range("C10").Value = "C1A3"
$newSheet.Range("D10").Value = "C1B7"
$newSheet.Range("E10:K10").ClearContents()

# --- Only 9 data columns are used this time (vs 10 previously): drop column L ---
$newSheet.Range("L8:L18").ClearContents()

# --- Row 18 (previously the "10" row marker + data) is now unused/blank ---
$newSheet.Range("B18:K18").ClearContents()

# --- Remove the borders from the now-empty column L and row 18 ---
$newSheet.Range("L8:L18").Borders.LineStyle = 0
$newSheet.Range("B18:L18").Borders.LineStyle = 0

# --- The "labeled tape" footer row no longer has a top border ---
$newSheet.Range("B19:L19").Borders.LineStyle = 0

# --- Selection / active-tab bookkeeping to mirror the edit ---
$oldSheet = $wb.Worksheets.Item("R3_P6")
$oldSheet.Activate()
$excel.ActiveWindow.ScrollRow = 2
$oldSheet.Range("D21").Select()

$newSheet.Activate()
$newSheet.Range("B6,J5").Select()
$newSheet.Range("J5").Activate()
